# Fix Training Data Issue (#48)
# Data was taken from 1 day off due to way NBA stats were shown.
# Column BF ("Date") on Sheet1 contained a malformed date string
# ("6-23-2013-14"); correct it to the proper ISO date "2014-06-23"
# for every data row (rows 2-31), keeping the cell a plain text value
# (not letting Excel auto-convert it into a date serial number) and
# without leaving any residual formatting/style on the cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 31
$col = 58  # column BF

$rng = $ws.Range($ws.Cells.Item($firstRow, $col), $ws.Cells.Item($lastRow, $col))

# Force text storage so "2014-06-23" isn't reinterpreted as a date serial.
$rng.NumberFormat = "@"

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, $col).Value = "2014-06-23"
}

# Restore the default "Normal" cell style so no stray formatting remains
# on the cells (they were unstyled before the edit).
$rng.Style = "Normal"
